$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet/tab
$ws.Name = "Users.xlsx"

# New rows (18-22) need to inherit the existing row styling (A/F numeric
# style, B/C/D/E text style) before we fill in their values.
$ws.Range("A17:F17").Copy()
$ws.Range("A18:F22").PasteSpecial(-4122)  # xlPasteFormats

# Data rows (row, FirstName, LastName, Email, Phone, Year) - rows 2 & 3 unchanged
$data = @(
    @(4,  "User5",       "User5LN",        "user5@gmail.com",       "(159)753852", 2000),
    @(5,  "User2",       "LastUser2",      "user2@gmail.com",       "123456789",   1923),
    @(6,  "User8",       "User8LN",        "user8@gmail.com",       "(159)753852", 2000),
    @(7,  "Boxis",       "Strong",         "boxis@gmail.com",       "(111)333222", 1993),
    @(8,  "Tor",         "Asgaard",        "tor@gmail.com",         "(111)333888", 1994),
    @(9,  "User1",       "Admin1",         "user1@gmail.com",       "(123)456780", 1990),
    @(10, "Gunnar",      "Jensen",         "gunnar@gmail.com",      "(111)222444", 1980),
    @(11, "TestTrainer", "TeatTrainerLN",  "testTrainer@gmail.com", "(123)123123", 1998),
    @(12, "Bruce",       "Lee",            "bruce@gmail.com",       "(111)333445", 1987),
    @(13, "Gamora",      "Gamorak",        "gamora@gmail.com",      "(111)333111", 1988),
    @(14, "Witcher",     "Moon",           "witcher@gmail.com",     "(111)333999", 1990),
    @(15, "Supwom",      "Nanual",         "supwom@gmail.com",      "(111)333777", 1988),
    @(16, "Barney",      "Ross",           "barney@gmail.com",      "(111)222333", 1975),
    @(17, "Lee",         "Christmas",      "lee@gmail.com",         "(111)333444", 1977),
    @(18, "Marvel",      "Levram",         "marvel@gmail.com",      "(111)333555", 1995),
    @(19, "Jean",        "Vilain",         "jean@gmail.com",        "(111)222777", 1973),
    @(20, "User6",       "User6LN",        "user6@gmail.com",       "",            0),
    @(21, "Sonya",       "Night",          "sonya@gmail.com",       "(111)333666", 1996),
    @(22, "Natalia",     "Romanoff",       "natalia@gmail.com",     "(111)222888", 1986)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $r - 1
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

# Row 5 (Phone column) holds a digits-only value that must stay text
# ("123456789"), not be auto-converted to a number. Format as text before
# assigning it, then re-copy the neighbouring cell's formatting back on
# top so the cell keeps the same shared text style as the rest of the
# column (NumberFormat="@" alone would allocate a distinct style).
$ws.Cells.Item(5, 5).NumberFormat = "@"
$ws.Cells.Item(5, 5).Value = "123456789"
$ws.Range("D5").Copy()
$ws.Range("E5").PasteSpecial(-4122)  # xlPasteFormats
